# Fruta / hortaliza, semanal
#
# Inserts a new weekly price record for Chirimoya ("Vega Monumental
# Concepción") at row 7, pushing the previously-existing rows 7..29 down
# to rows 8..30 (dimension grows from A1:T29 to A1:T30).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 7..29 down to 8..30 and leave a fresh, blank row 7 for the
# new record. Excel carries the existing cell formatting (e.g. the date
# number format on column D) down with the shifted rows.
$ws.Rows.Item(7).Insert()

# Populate the newly inserted row 7 with the new weekly record.
$ws.Cells.Item(7, 1).Value = 11
$ws.Cells.Item(7, 2).Value = 'Vega Monumental Concepción'
$ws.Cells.Item(7, 3).Value = 'Bíobío'
$ws.Cells.Item(7, 4).Value = 44819
$ws.Cells.Item(7, 5).Value = 8
$ws.Cells.Item(7, 6).Value = 'Fruta'
$ws.Cells.Item(7, 7).Value = 100107
$ws.Cells.Item(7, 8).Value = 'Otros'
$ws.Cells.Item(7, 9).Value = 100107002
$ws.Cells.Item(7, 10).Value = 'Chirimoya'
$ws.Cells.Item(7, 11).Value = 'Cultivar IV Región'
$ws.Cells.Item(7, 12).Value = 'Primera'
$ws.Cells.Item(7, 13).Value = 80
$ws.Cells.Item(7, 14).Value = 25000
$ws.Cells.Item(7, 15).Value = 26000
$ws.Cells.Item(7, 16).Value = 25500
$ws.Cells.Item(7, 17).Value = '$/bandeja 10 kilos'
$ws.Cells.Item(7, 18).Value = 'Provincia de Limarí'
$ws.Cells.Item(7, 19).Value = 2550
$ws.Cells.Item(7, 20).Value = 10
